$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fields = @(
    @{Col='S'; Header='biotic_relationship'; Comment='Free-living or from host (define relationship)'},
    @{Col='T'; Header='chem_administration'; Comment='list of chemical compounds administered to the host or site where sampling occurred, and when (e.g. antibiotics, N fertilizer, air filter); can include multiple compounds. For Chemical Entities of Biological Interest ontology (CHEBI) (v1.72), please see http://bioportal.bioontology.org/visualize/44603'},
    @{Col='U'; Header='dermatology_disord'; Comment='history of dermatology disorders; can include multiple disorders'},
    @{Col='V'; Header='dominant_hand'; Comment='dominant hand of the subject'},
    @{Col='W'; Header='ethnicity'; Comment='ethnicity of the subject'},
    @{Col='X'; Header='extrachrom_elements'; Comment='Plasmids that have significance phenotypic consequence'},
    @{Col='Y'; Header='host_age'; Comment='Age of host at the time of sampling'},
    @{Col='Z'; Header='host_body_mass_index'; Comment='body mass index of the host, calculated as weight/(height)squared'},
    @{Col='AA'; Header='host_body_product'; Comment='substance produced by the host, e.g. stool, mucus, where the sample was obtained from'},
    @{Col='AB'; Header='host_body_temp'; Comment='core body temperature of the host when sample was collected'},
    @{Col='AC'; Header='host_diet'; Comment='type of diet depending on the sample for animals omnivore, herbivore etc., for humans high-fat, meditteranean etc.; can include multiple diet types'},
    @{Col='AD'; Header='host_disease'; Comment='Name of relevant disease, e.g. Salmonella gastroenteritis. For the controlled vocabulary, please see Human Disease Ontology, http://bioportal.bioontology.org/ontologies/1009 or MeSH, http://www.ncbi.nlm.nih.gov/mesh'},
    @{Col='AE'; Header='host_family_relationship'; Comment=$null},
    @{Col='AF'; Header='host_genotype'; Comment=$null},
    @{Col='AG'; Header='host_height'; Comment='the height of subject'},
    @{Col='AH'; Header='host_last_meal'; Comment='content of last meal and time since feeding; can include multiple values'},
    @{Col='AI'; Header='host_occupation'; Comment='most frequent job performed by subject'},
    @{Col='AJ'; Header='host_phenotype'; Comment=$null},
    @{Col='AK'; Header='host_pulse'; Comment='resting pulse of the host, measured as beats per minute'},
    @{Col='AL'; Header='host_sex'; Comment='Gender or physical sex of the host'},
    @{Col='AM'; Header='host_subject_id'; Comment='a unique identifier by which each subject can be referred to, de-identified, e.g. #131'},
    @{Col='AN'; Header='host_tissue_sampled'; Comment='Type of tissue the initial sample was taken from. Controlled vocabulary, http://bioportal.bioontology.org/ontologies/1005'},
    @{Col='AO'; Header='host_tot_mass'; Comment='total mass of the host at collection, the unit depends on host'},
    @{Col='AP'; Header='ihmc_medication_code'; Comment='can include multiple medication codes'},
    @{Col='AQ'; Header='isolation_source'; Comment='Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.'},
    @{Col='AR'; Header='medic_hist_perform'; Comment='whether full medical history was collected'},
    @{Col='AS'; Header='misc_param'; Comment='any other measurement performed or parameter collected, that is not listed here'},
    @{Col='AT'; Header='organism_count'; Comment='total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts'},
    @{Col='AU'; Header='oxy_stat_samp'; Comment='oxygenation status of sample'},
    @{Col='AV'; Header='perturbation'; Comment='type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types'},
    @{Col='AW'; Header='rel_to_oxygen'; Comment='Aerobic or anaerobic'},
    @{Col='AX'; Header='samp_collect_device'; Comment='Method or device employed for collecting sample'},
    @{Col='AY'; Header='samp_mat_process'; Comment='Processing applied to the sample during or after isolation'},
    @{Col='AZ'; Header='samp_salinity'; Comment='salinity of sample, i.e. measure of total salt concentration'},
    @{Col='BA'; Header='samp_size'; Comment='Amount or size of sample (volume, mass or area) that was collected'},
    @{Col='BB'; Header='samp_store_dur'; Comment='duration for which sample was stored'},
    @{Col='BC'; Header='samp_store_loc'; Comment='location at which sample was stored, usually name of a specific freezer/room'},
    @{Col='BD'; Header='samp_store_temp'; Comment='temperature at which sample was stored, e.g. -80'},
    @{Col='BE'; Header='samp_vol_we_dna_ext'; Comment='volume (mL) or weight (g) of sample processed for DNA extraction'},
    @{Col='BF'; Header='source_material_id'; Comment='unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.'},
    @{Col='BG'; Header='subspecf_gen_lin'; Comment='Information about the genetic distinctness of the lineage (eg., biovar, serovar)'},
    @{Col='BH'; Header='temperature'; Comment='temperature of the sample at time of sampling'},
    @{Col='BI'; Header='time_since_last_wash'; Comment='specification of the time since last wash'},
    @{Col='BJ'; Header='trophic_level'; Comment='Feeding position in food chain (eg., chemolithotroph)'}
)

foreach ($f in $fields) {
    $addr = $f.Col + "15"
    $ws.Range($addr).Value = $f.Header
    if ($f.Comment -ne $null) {
        $ws.Range($addr).AddComment($f.Comment)
    }
}

# Copy the existing "optional field" (yellow) format onto the newly added header cells
$ws.Range("C15").Copy()
$ws.Range("S15:BJ15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
